# Update countries & provincias Spain
# - Refresh the "updated at" timestamp in A1.
# - Refresh COVID-19 daily figures for several countries (rows keep their
#   sort order by "Casos totales"; the countries occupying rows 77-80
#   change because Honduras' case count overtook Croacia / Bosnia y
#   Herzegovina / Senegal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 04:05"

# --- Panama (row 50) ---------------------------------------------------
$ws.Range("B50").Value = 8944
$ws.Range("C50").Value = 161
$ws.Range("D50").Value = 6067
$ws.Range("E50").Value = 2621
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 256

# --- Rows 77-80: Honduras climbs above Croacia / Bosnia / Senegal -----
# Row 77 now holds Honduras (previously Croacia)
$ws.Range("A77").Value = "Honduras"
$ws.Range("B77").Value = 2255
$ws.Range("C77").Value = 175
$ws.Range("D77").Value = 237
$ws.Range("E77").Value = 1895
$ws.Range("F77").Value = 10
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 123

# Row 78 now holds Croacia (previously Bosnia y Herzegovina)
$ws.Range("A78").Value = "Croacia"
$ws.Range("B78").Value = 2213
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 1834
$ws.Range("E78").Value = 285
$ws.Range("F78").Value = 9
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 94

# Row 79 now holds Bosnia y Herzegovina (previously Senegal)
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 2181
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 1228
$ws.Range("E79").Value = 833
$ws.Range("F79").Value = 4
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 120

# Row 80 now holds Senegal (previously Honduras)
$ws.Range("A80").Value = "Senegal"
$ws.Range("B80").Value = 2105
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 782
$ws.Range("E80").Value = 1302
$ws.Range("F80").Value = 6
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 21

# --- Isla de Man (row 134) ---------------------------------------------
$ws.Range("D134").Value = 274
$ws.Range("E134").Value = 35

# --- Nepal (row 142) ----------------------------------------------------
$ws.Range("B142").Value = 245
$ws.Range("C142").Value = 2
$ws.Range("E142").Value = 210
